$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# The Status text for the 434a3989 file changes from "Ready for handoff" to
# "Handback transform failed" everywhere it's shown (Overview B3/C3, and the
# Status column (C3) on both the zh-cn and de-de detail sheets).
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# New "Error Detail" (column K) entries on row 3 of each locale sheet,
# describing the handback/handoff file name mismatch for that locale.
$zhcn.Range("K3").Value = "Handback file name: vu1sb2ju.q3g is different with handoff file name: 434a3989-d849-4143-bfd2-c0d3718e0e20.bad8a287817f5aed8e48bafb3f793fb9b6309caf.zh-cn."
$dede.Range("K3").Value = "Handback file name: vu1sb2ju.q3g is different with handoff file name: 434a3989-d849-4143-bfd2-c0d3718e0e20.bad8a287817f5aed8e48bafb3f793fb9b6309caf.de-de."
